$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.782.95"
$ws.Cells.Item(2, 5).Value = "  +4.76%  "

$ws.Cells.Item(3, 4).Value = "1.615.00"
$ws.Cells.Item(3, 5).Value = "  +4.26%  "

$ws.Cells.Item(4, 5).Value = "  -0.46%  "

$ws.Cells.Item(5, 4).Value = "214.68"
$ws.Cells.Item(5, 5).Value = "  +2.01%  "

$ws.Cells.Item(6, 5).Value = "  +7.38%  "

$ws.Cells.Item(7, 4).Value = "0.996"
$ws.Cells.Item(7, 5).Value = "  -0.41%  "

$ws.Cells.Item(8, 5).Value = "  +12.89%  "

$ws.Cells.Item(9, 5).Value = "  +3.75%  "

$ws.Cells.Item(10, 5).Value = "  +3.04%  "

$ws.Cells.Item(11, 5).Value = "  +2.97%  "

$ws.Cells.Item(12, 4).Value = "1.850.30"
$ws.Cells.Item(12, 5).Value = "  +4.52%  "

$ws.Cells.Item(13, 4).Value = "1.612.08"
$ws.Cells.Item(13, 5).Value = "  +3.73%  "

$ws.Cells.Item(14, 4).Value = "29.787.51"
$ws.Cells.Item(14, 5).Value = "  +4.89%  "

$ws.Cells.Item(15, 4).Value = "0.539"
$ws.Cells.Item(15, 5).Value = "  +6.18%  "

$ws.Cells.Item(16, 4).Value = "3.77"
$ws.Cells.Item(16, 5).Value = "  +4.04%  "

$ws.Cells.Item(17, 4).Value = "246.24"
$ws.Cells.Item(17, 5).Value = "  +7.92%  "

$ws.Cells.Item(18, 4).Value = "63.72"
$ws.Cells.Item(18, 5).Value = "  +4.57%  "

$ws.Cells.Item(19, 4).Value = "7.64"
$ws.Cells.Item(19, 5).Value = "  +4.54%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0697"
$ws.Cells.Item(20, 5).Value = "  +3.84%  "

$ws.Cells.Item(21, 4).Value = "0.993"
$ws.Cells.Item(21, 5).Value = "  -0.75%  "

$ws.Cells.Item(22, 4).Value = "4.07"
$ws.Cells.Item(22, 5).Value = "  +4.78%  "

$ws.Cells.Item(23, 4).Value = "9.29"
$ws.Cells.Item(23, 5).Value = "  +4.31%  "

$ws.Cells.Item(24, 5).Value = "  +5.19%  "

$ws.Cells.Item(25, 4).Value = "155.98"
$ws.Cells.Item(25, 5).Value = "  +2.87%  "

$ws.Cells.Item(26, 4).Value = "15.39"
$ws.Cells.Item(26, 5).Value = "  +4.48%  "

$ws.Cells.Item(27, 5).Value = "  +6.15%  "

$ws.Cells.Item(28, 5).Value = "  +3.53%  "

$ws.Cells.Item(29, 4).Value = "0.997"
$ws.Cells.Item(29, 5).Value = "  -0.39%  "

$ws.Cells.Item(30, 4).Value = "0.0474"
$ws.Cells.Item(30, 5).Value = "  +1.69%  "

$ws.Cells.Item(31, 5).Value = "  +0.56%  "

$ws.Cells.Item(32, 5).Value = "  +3.47%  "

$ws.Cells.Item(33, 4).Value = "1.446.11"
$ws.Cells.Item(33, 5).Value = "  +4.53%  "

$ws.Cells.Item(34, 4).Value = "3.12"
$ws.Cells.Item(34, 5).Value = "  +4.22%  "

$ws.Cells.Item(35, 5).Value = "  -0.64%  "

$ws.Cells.Item(36, 2).Value = "MXToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(36, 4).Value = "2.85"
$ws.Cells.Item(36, 5).Value = "  +11.01%  "

$ws.Cells.Item(37, 2).Value = "LidoDAOToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(37, 4).Value = "1.52"
$ws.Cells.Item(37, 5).Value = "  +3.50%  "

$ws.Cells.Item(38, 4).Value = "2.32"
$ws.Cells.Item(38, 5).Value = "  +0.68%  "

$ws.Cells.Item(39, 5).Value = "  +3.61%  "

$ws.Cells.Item(40, 4).Value = "56.59"

$ws.Cells.Item(41, 4).Value = "0.539"
$ws.Cells.Item(41, 5).Value = "  +5.81%  "

$ws.Cells.Item(42, 4).Value = "1.96"
$ws.Cells.Item(42, 5).Value = "  +2.00%  "

$ws.Cells.Item(43, 4).Value = "68.88"
$ws.Cells.Item(43, 5).Value = "  +11.53%  "

$ws.Cells.Item(44, 4).Value = "0.801"
$ws.Cells.Item(44, 5).Value = "  +4.08%  "

$ws.Cells.Item(45, 4).Value = "0.996"
$ws.Cells.Item(45, 5).Value = "  -0.38%  "

$ws.Cells.Item(46, 5).Value = "  +2.53%  "

$ws.Cells.Item(47, 4).Value = "5.36"
$ws.Cells.Item(47, 5).Value = "  +0.85%  "

$ws.Cells.Item(48, 4).Value = "1.757.57"
$ws.Cells.Item(48, 5).Value = "  +4.45%  "

$ws.Cells.Item(49, 4).Value = "86.80"
$ws.Cells.Item(49, 5).Value = "  +1.56%  "

$ws.Cells.Item(50, 4).Value = "0.838"
$ws.Cells.Item(50, 5).Value = "  -4.06%  "

$ws.Cells.Item(51, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(51, 4).Value = "0.0₆0103"
$ws.Cells.Item(51, 5).Value = "  +2.54%  "
